$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.789.00'
$ws.Range("E2").Value = '  +1.68%  '
$ws.Range("D3").Value = '2.801.47'
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'350.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = "'112.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.14%  '
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = "'0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.17%  '
$ws.Range("D10").Value = "'40.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.43%  '
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("E12").Value = '  +0.68%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").Value = "'7.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.74%  '
$ws.Range("D15").Value = '3.237.38'
$ws.Range("E15").Value = '  +1.59%  '
$ws.Range("D16").Value = "'0.968"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.26%  '
$ws.Range("D17").Value = '2.782.36'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '51.770.86'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = "'3.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.70%  '
$ws.Range("D20").Value = "'7.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").Value = "'13.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.95%  '
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D24").Value = "'269.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.84%  '
$ws.Range("D25").Value = "'2.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.04%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  +13.78%  '
$ws.Range("D30").Value = "'10.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.77%  '
$ws.Range("D31").Value = "'2.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("E32").Value = '  +2.10%  '
$ws.Range("D33").Value = "'6.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("D34").Value = "'0.0910"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.82%  '
$ws.Range("D35").Value = "'0.0454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("E36").Value = '  +5.29%  '
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'18.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.55%  '
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("D42").Value = "'2.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("D43").Value = "'122.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").Value = "'2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.73%  '
$ws.Range("D45").Value = "'22.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("D46").Value = "'3.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.67%  '
$ws.Range("D47").Value = "'2.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.81%  '
$ws.Range("D48").Value = '2.121.12'
$ws.Range("D49").Value = "'0.986"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.50%  '
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("E51").Value = '  +17.31%  '
